$wb = $excel.ActiveWorkbook

# --- Spirits sheet: append new rows (Name column) ---
$ws = $wb.Worksheets.Item("Spirits")
$spirits = @("Smirnoff","Gordons","Hendricks","Bacardi","Jameson","Powers","Paddy","Hennessy","Huzzar")
for ($i = 0; $i -lt $spirits.Length; $i++) {
    $ws.Range("A" + (2 + $i)).Value = $spirits[$i]
}
$ws.Activate()
$ws.Range("A11").Select()

# --- Beer sheet: append new rows (Name column) ---
$ws = $wb.Worksheets.Item("Beer")
$beers = @("Guinness","Coors","Heineken","Smithwicks","Harp","Murphys","Beamish","Carlsberg")
for ($i = 0; $i -lt $beers.Length; $i++) {
    $ws.Range("A" + (2 + $i)).Value = $beers[$i]
}
$ws.Range("A10").Select()

# --- Misc sheet: append new rows (Name column) ---
$ws = $wb.Worksheets.Item("Misc")
$misc = @("Coca Cola","Fanta Orange","Fanta Lemon","Club Lemon","Club Orange","7 Up","Sprite","Tonic water","Schweppes Tonic","Schweppes Slimline Tonic","Slimline tonic","Ginger ale","Lucozade","Red Bull")
for ($i = 0; $i -lt $misc.Length; $i++) {
    $ws.Range("A" + (2 + $i)).Value = $misc[$i]
}
$ws.Range("A16").Select()

# --- Heroes sheet: insert new "Icon" column before the Image Filename column (G) ---
$ws = $wb.Worksheets.Item("Heroes")
$ws.Columns("G").Insert()
$ws.Range("G1").Value = "Icon"
$ws.Columns("G").ColumnWidth = 29.7
$ws.Activate()
$ws.Range("W1").Select()

# --- Cocktails sheet: selection change only ---
$ws = $wb.Worksheets.Item("Cocktails")
$ws.Range("B1").Select()

# --- Make Spirits the active sheet/tab (matches activeTab change) ---
$ws = $wb.Worksheets.Item("Spirits")
$ws.Activate()
